$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Fgf2-Fgfr1 LR-pair results for YoungD0, following Dr Hou advice:
# the table now reports the full 3x3 sending x target cluster combinations
# (ECs, FAPs, sCs) instead of a subset, with refreshed statistics.
# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf2"
$ws.Cells.Item(2,3).Value = "Fgfr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.6462393333333333
$ws.Cells.Item(2,8).Value = 1.938718
$ws.Cells.Item(2,9).Value = 0.03461850536298827
$ws.Cells.Item(2,10).Value = 0.03461850536298827
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.675378666666666
$ws.Cells.Item(2,14).Value = 14.026136
$ws.Cells.Item(2,15).Value = 0.03681964474327726
$ws.Cells.Item(2,16).Value = 0.03681964474327726
$ws.Cells.Item(2,17).Value = 3.021413592627555
$ws.Cells.Item(2,18).Value = 27.192722333648
$ws.Cells.Item(2,19).Value = 0.001274641069008467
$ws.Cells.Item(2,20).Value = 0.001274641069008467

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf2"
$ws.Cells.Item(3,3).Value = "Fgfr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.6462393333333333
$ws.Cells.Item(3,8).Value = 1.938718
$ws.Cells.Item(3,9).Value = 0.03461850536298827
$ws.Cells.Item(3,10).Value = 0.03461850536298827
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 82.95722966666666
$ws.Cells.Item(3,14).Value = 248.871689
$ws.Cells.Item(3,15).Value = 0.653306596744776
$ws.Cells.Item(3,16).Value = 0.653306596744776
$ws.Cells.Item(3,17).Value = 53.61022479496688
$ws.Cells.Item(3,18).Value = 482.492023154702
$ws.Cells.Item(3,19).Value = 0.02261649792308464
$ws.Cells.Item(3,20).Value = 0.02261649792308464

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf2"
$ws.Cells.Item(4,3).Value = "Fgfr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.6462393333333333
$ws.Cells.Item(4,8).Value = 1.938718
$ws.Cells.Item(4,9).Value = 0.03461850536298827
$ws.Cells.Item(4,10).Value = 0.03461850536298827
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 39.34793966666667
$ws.Cells.Item(4,14).Value = 118.043819
$ws.Cells.Item(4,15).Value = 0.3098737585119468
$ws.Cells.Item(4,16).Value = 0.3098737585119468
$ws.Cells.Item(4,17).Value = 25.42818629822689
$ws.Cells.Item(4,18).Value = 228.853676684042
$ws.Cells.Item(4,19).Value = 0.01072736637089516
$ws.Cells.Item(4,20).Value = 0.01072736637089516

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf2"
$ws.Cells.Item(5,3).Value = "Fgfr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 15.322826
$ws.Cells.Item(5,8).Value = 45.968478
$ws.Cells.Item(5,9).Value = 0.8208310864042159
$ws.Cells.Item(5,10).Value = 0.8208310864042158
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.675378666666666
$ws.Cells.Item(5,14).Value = 14.026136
$ws.Cells.Item(5,15).Value = 0.03681964474327726
$ws.Cells.Item(5,16).Value = 0.03681964474327726
$ws.Cells.Item(5,17).Value = 71.64001379344533
$ws.Cells.Item(5,18).Value = 644.7601241410081
$ws.Cells.Item(5,19).Value = 0.03022270899564155
$ws.Cells.Item(5,20).Value = 0.03022270899564155

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf2"
$ws.Cells.Item(6,3).Value = "Fgfr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 15.322826
$ws.Cells.Item(6,8).Value = 45.968478
$ws.Cells.Item(6,9).Value = 0.8208310864042159
$ws.Cells.Item(6,10).Value = 0.8208310864042158
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 82.95722966666666
$ws.Cells.Item(6,14).Value = 248.871689
$ws.Cells.Item(6,15).Value = 0.653306596744776
$ws.Cells.Item(6,16).Value = 0.653306596744776
$ws.Cells.Item(6,17).Value = 1271.139195624371
$ws.Cells.Item(6,18).Value = 11440.25276061934
$ws.Cells.Item(6,19).Value = 0.5362543635610554
$ws.Cells.Item(6,20).Value = 0.5362543635610554

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf2"
$ws.Cells.Item(7,3).Value = "Fgfr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 15.322826
$ws.Cells.Item(7,8).Value = 45.968478
$ws.Cells.Item(7,9).Value = 0.8208310864042159
$ws.Cells.Item(7,10).Value = 0.8208310864042158
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 39.34793966666667
$ws.Cells.Item(7,14).Value = 118.043819
$ws.Cells.Item(7,15).Value = 0.3098737585119468
$ws.Cells.Item(7,16).Value = 0.3098737585119468
$ws.Cells.Item(7,17).Value = 602.9216329708314
$ws.Cells.Item(7,18).Value = 5426.294696737483
$ws.Cells.Item(7,19).Value = 0.2543540138475189
$ws.Cells.Item(7,20).Value = 0.2543540138475189

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fgf2"
$ws.Cells.Item(8,3).Value = "Fgfr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.698388
$ws.Cells.Item(8,8).Value = 8.095164
$ws.Cells.Item(8,9).Value = 0.1445504082327959
$ws.Cells.Item(8,10).Value = 0.1445504082327959
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.675378666666666
$ws.Cells.Item(8,14).Value = 14.026136
$ws.Cells.Item(8,15).Value = 0.03681964474327726
$ws.Cells.Item(8,16).Value = 0.03681964474327726
$ws.Cells.Item(8,17).Value = 12.61598568958933
$ws.Cells.Item(8,18).Value = 113.543871206304
$ws.Cells.Item(8,19).Value = 0.005322294678627246
$ws.Cells.Item(8,20).Value = 0.005322294678627246

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fgf2"
$ws.Cells.Item(9,3).Value = "Fgfr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.698388
$ws.Cells.Item(9,8).Value = 8.095164
$ws.Cells.Item(9,9).Value = 0.1445504082327959
$ws.Cells.Item(9,10).Value = 0.1445504082327959
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 82.95722966666666
$ws.Cells.Item(9,14).Value = 248.871689
$ws.Cells.Item(9,15).Value = 0.653306596744776
$ws.Cells.Item(9,16).Value = 0.653306596744776
$ws.Cells.Item(9,17).Value = 223.8507930457773
$ws.Cells.Item(9,18).Value = 2014.657137411996
$ws.Cells.Item(9,19).Value = 0.09443573526063592
$ws.Cells.Item(9,20).Value = 0.09443573526063592

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fgf2"
$ws.Cells.Item(10,3).Value = "Fgfr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.698388
$ws.Cells.Item(10,8).Value = 8.095164
$ws.Cells.Item(10,9).Value = 0.1445504082327959
$ws.Cells.Item(10,10).Value = 0.1445504082327959
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 39.34793966666667
$ws.Cells.Item(10,14).Value = 118.043819
$ws.Cells.Item(10,15).Value = 0.3098737585119468
$ws.Cells.Item(10,16).Value = 0.3098737585119468
$ws.Cells.Item(10,17).Value = 106.1760082212573
$ws.Cells.Item(10,18).Value = 955.584073991316
$ws.Cells.Item(10,19).Value = 0.04479237829353272
$ws.Cells.Item(10,20).Value = 0.04479237829353271
